# The deck's Notes Master has a "Date" placeholder containing an
# auto-updating date field (type="datetimeFigureOut"). Its cached display
# value needs to change from "5/12/25" to "7/14/25" (e.g. the file was
# re-saved/opened on a later date, refreshing the cached text).
#
# The canonical way to control this placeholder's text through the
# PowerPoint object model is via HeadersFooters.DateAndTime on the Notes
# Master (Insert > Header & Footer > Notes and Handouts > Date and time).

$p  = $ppt.ActivePresentation
$nm = $p.NotesMaster

$nm.HeadersFooters.DateAndTime.Text = "7/14/25"
